# Artfynd observation export: the underlying report rows (identified by the
# "Id" value in column A) were re-shuffled into different row positions by
# the source system. Re-apply that shuffle by writing each destination row's
# final field values directly.
#
# For text-typed cells whose value looks like a plain number (column I -
# "Antal"/count - and the empty placeholders that replace it), a leading
# apostrophe is used so Excel stores the cell as Text (matching the original
# t="inlineStr" cell type) instead of re-interpreting it as a Number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 9's data: Id 111618046, Blåmossa)
$ws.Range("A2").Value = 111618046
$ws.Range("B2").Value = 93388
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 2180
$ws.Range("F2").Value = "Blåmossa"
$ws.Range("G2").Value = "Leucobryum glaucum"
$ws.Range("H2").Value = "(Hedw.) Ångstr."
$ws.Range("I2").Value = "'"
$ws.Range("J2").Value = "'"
$ws.Range("K2").Value = "'"
$ws.Range("Q2").Value = 580591.6383206119
$ws.Range("R2").Value = 6415156.322361182
$ws.Range("AC2").Value = "'"

# Row 3 (was row 2's data: Id 111618089, Knärot)
$ws.Range("A3").Value = 111618089
$ws.Range("I3").Value = "'30"
$ws.Range("P3").Value = "A 32649, Heda, Sm"
$ws.Range("Q3").Value = 580617.6201989455
$ws.Range("R3").Value = 6415136.627037819

# Row 4 (was row 5's data: Id 111618039, Blåmossa)
$ws.Range("A4").Value = 111618039
$ws.Range("B4").Value = 93388
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 2180
$ws.Range("F4").Value = "Blåmossa"
$ws.Range("G4").Value = "Leucobryum glaucum"
$ws.Range("H4").Value = "(Hedw.) Ångstr."
$ws.Range("J4").Value = "'"
$ws.Range("Q4").Value = 580599.6803078586
$ws.Range("R4").Value = 6415233.627682217

# Row 5 (was row 4's data: Id 111618078, Knärot)
$ws.Range("A5").Value = 111618078
$ws.Range("B5").Value = 96348
$ws.Range("D5").Value = "VU"
$ws.Range("E5").Value = 220787
$ws.Range("F5").Value = "Knärot"
$ws.Range("G5").Value = "Goodyera repens"
$ws.Range("H5").Value = "(L.) R. Br."
$ws.Range("J5").Value = "plantor/tuvor"
$ws.Range("Q5").Value = 580612.1009209087
$ws.Range("R5").Value = 6415119.491031807

# Row 6 (was row 7's data: Id 111618070, Knärot)
$ws.Range("A6").Value = 111618070
$ws.Range("I6").Value = "'15"
$ws.Range("K6").Value = "blomning"
$ws.Range("Q6").Value = 580592.470229132
$ws.Range("R6").Value = 6415141.442167919
$ws.Range("AC6").Value = "1 blomma"

# Row 7 (was row 8's data: Id 111618056, Knärot)
$ws.Range("A7").Value = 111618056
$ws.Range("Q7").Value = 580582.6881743574
$ws.Range("R7").Value = 6415124.22061418
$ws.Range("AC7").Value = "2 blommor"

# Row 8 (was row 3's data: Id 111618109, Knärot)
$ws.Range("A8").Value = 111618109
$ws.Range("I8").Value = "'10"
$ws.Range("P8").Value = "A 32649, Sm"
$ws.Range("Q8").Value = 580619.1666838422
$ws.Range("R8").Value = 6415112.716507593
$ws.Range("AC8").Value = "1 blomma"

# Row 9 (was row 6's data: Id 111618144, Knärot)
$ws.Range("A9").Value = 111618144
$ws.Range("B9").Value = 96348
$ws.Range("D9").Value = "VU"
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = "Knärot"
$ws.Range("G9").Value = "Goodyera repens"
$ws.Range("H9").Value = "(L.) R. Br."
$ws.Range("I9").Value = "'2"
$ws.Range("J9").Value = "plantor/tuvor"
$ws.Range("Q9").Value = 580620.6996611424
$ws.Range("R9").Value = 6415142.541277731
